$d = $word.ActiveDocument

$replacements = @(
    @("820÷4=205, 0", "579÷7=82, 5"),
    @("740÷9=82, 2", "690÷8=86, 2"),
    @("176÷8=22, 0", "766÷4=191, 2"),
    @("412÷7=58, 6", "937÷3=312, 1"),
    @("644÷9=71, 5", "748÷8=93, 4"),
    @("863÷8=107, 7", "688÷4=172, 0"),
    @("795÷7=113, 4", "263÷5=52, 3"),
    @("601÷6=100, 1", "675÷7=96, 3"),
    @("546÷3=182, 0", "554÷7=79, 1"),
    @("264÷8=33, 0", "583÷7=83, 2"),
    @("534÷8=66, 6", "586÷4=146, 2"),
    @("291÷3=97, 0", "681÷8=85, 1"),
    @("935÷4=233, 3", "396÷6=66, 0"),
    @("544÷6=90, 4", "542÷4=135, 2"),
    @("278÷6=46, 2", "559÷4=139, 3"),
    @("650÷9=72, 2", "332÷5=66, 2"),
    @("345÷9=38, 3", "768÷6=128, 0"),
    @("784÷8=98, 0", "175÷7=25, 0"),
    @("513÷7=73, 2", "140÷7=20, 0"),
    @("499÷3=166, 1", "496÷2=248, 0"),
    @("951÷2=475, 1", "976÷9=108, 4"),
    @("133÷8=16, 5", "567÷4=141, 3"),
    @("683÷8=85, 3", "775÷3=258, 1"),
    @("267÷9=29, 6", "644÷5=128, 4"),
    @("623÷2=311, 1", "339÷8=42, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
